$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the results for row 12 (final event - Spartan challenge)
$ws.Range("B12").Value = 375
$ws.Range("C12").Value = 355
$ws.Range("D12").Value = 380

$ws.Range("H12").Value = 262
$ws.Range("I12").Value = 339
$ws.Range("J12").Value = 339

$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("R12").Value = 0
$ws.Range("S12").Value = 0

$ws.Range("T12").Value = 214
$ws.Range("U12").Value = 204
$ws.Range("V12").Value = 234

$ws.Range("W12").Value = 100
$ws.Range("X12").Value = 150
$ws.Range("Y12").Value = 175
$ws.Range("Z12").Value = 0
$ws.Range("AA12").Value = 0
$ws.Range("AB12").Value = 0
$ws.Range("AC12").Value = 200

# Update the active cell selection to A12
$ws.Range("A12").Select()
